# Author's edit: after the paragraph that ends in "...avaliação das
# métricas de performance." (the last sentence of the "Para a partição
# do conjunto de dados..." paragraph), a new, completely empty paragraph
# is inserted, using the same "standard body" paragraph formatting used
# throughout this section (autoSpaceDE/DN off, adjustRightInd off,
# spacing after 0 / line 240 auto, justified).

$d = $word.ActiveDocument

$needle = "avaliação das métricas de performance."
$rng = $d.Content
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target sentence to anchor the new paragraph."
}

# Build a fresh Range object collapsed right after the matched text so we
# don't inherit any stray "current formatting" state left over from Find.
$insPoint = $d.Range($rng.End, $rng.End)

$newParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
              "<w:pPr>" + `
              "<w:autoSpaceDE w:val='0'/>" + `
              "<w:autoSpaceDN w:val='0'/>" + `
              "<w:adjustRightInd w:val='0'/>" + `
              "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/>" + `
              "<w:jc w:val='both'/>" + `
              "</w:pPr>" + `
              "</w:p>"

[void]$insPoint.InsertXML($newParaXml)

Write-Host "Inserted new empty paragraph after the sentence ending in '$needle'."
